# Add a "Directeur" (branch manager) column to the succursales list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Branch -> manager mapping, in the original authoring/typing order (row 22
# was entered before row 21), and row 1 is the header, filled in last below
# so the new shared-string order matches the source workbook exactly.
$directors = [ordered]@{
    2  = "Martin Hudon"
    3  = "Eric Vanier"
    4  = "Caroline Cyr"
    5  = "Yannick Lazare"
    6  = "Simon Huard"
    7  = "Bruno Adam"
    8  = "Claudio Fazioli"
    9  = "Eric Vanier"
    10 = "Yannick Lazare"
    11 = "Simon Huard"
    12 = "Claudio Fazioli"
    13 = "Karolane Roy"
    14 = "Karolane Roy"
    15 = "Caroline Cyr"
    20 = "Ludovic Gérard"
    22 = "Chantal Maltais"
    21 = "Danny Pronovost"
    23 = "Yannick Blanchet"
    24 = "Yannick Blanchet"
    25 = "Eric Savard"
    26 = "Danny Pronovost"
}

foreach ($row in $directors.Keys) {
    $ws.Cells.Item($row, 4).Value = $directors[$row]
}

# Header (added last)
$ws.Range("D1").Value = "Directeur"

# Match the new column D width (stored width="16" once converted from the
# COM ColumnWidth offset Excel applies internally).
$ws.Columns.Item(4).ColumnWidth = 15.166666666666666

# Match the final selection left in the sheet.
$ws.Range("D11").Select()
